# Timetracking sheet update: add a "Setup webappa" task (row 7, previously a
# blank gap row), record the extra time spent on the "create account" page
# (row 10) and mark both tasks with a footnote about bad working conditions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the previously-empty row 7 with new task data / extend row 10 data,
# entering values in the same order the new shared strings appear so the
# resulting string table ordering matches.
$ws.Range("C7").Value = "1h"
$ws.Range("D7").Value = "1h 10min"
$ws.Range("C10").Value = "1h 30min"
$ws.Range("D10").Value = "1h 50min"

# Update B10 text (append footnote marker)
$ws.Range("B10").Value = "Frontend - Stranica za create account *"

# New task name for row 7
$ws.Range("B7").Value = "Setup webappa *"

# Add footnote in F3
$ws.Range("F3").Value = "*rad u losim uslovima"

# Update selection / scroll position to match the saved view state
$ws.Range("E4").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
